$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph 1: turn the original (bookmark-interrupted, multi-run) paragraph
# into "2)" + the same sentence as a single clean run, with no bookmark.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1FullText = 'Considerando a necessidade do cliente, o precisa de um recurso que ele possa implementar rapidamente, o modelo que mais encaixa é o de prototipação evolucionaria, com uma entrega rápida de um protótipo, o cliente poderá implementar em um curto espaço de tempo um recurso na empresa, a ideia deste modelo é criar um protótipo e o projeto evoluir a partir do mesmo, podendo ser descartado, revisado ou aproveitado até a construção definitiva, possibilitando uma revisão e adaptação evolutiva, neste ciclo o desenvolvimento vai ocorrendo em ciclos e partes do produto podem ser desenvolvidas em separado e depois integradas, cobrindo outra necessidade do cliente, que precisa de diversos recursos na empresa.'

$p1Body = $d.Range($p1.Range.Start, $p1.Range.End - 1)
# Swap through a placeholder first: an identical-text assignment is a no-op,
# so this forces Word to rebuild the run (merging every fragment into one)
# and drops the _GoBack bookmark that used to sit inside the old runs.
$p1Body.Text = '.'
$p1Body2 = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$p1Body2.Text = $p1FullText

# Now prefix with its own "2)" run.
$p1Start = $d.Range($p1.Range.Start, $p1.Range.Start)
$p1Start.InsertBefore('2)')

# ---------------------------------------------------------------------------
# After paragraph 1: one blank paragraph, then the new "3)" paragraph (with
# the _GoBack bookmark moved into it), then two more blank paragraphs ahead
# of the pre-existing trailing blank paragraph.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$insertPoint = $d.Range($p1.Range.End, $p1.Range.End)

$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$p3Xml = "<w:p $wns>" +
         "<w:r><w:t>3)</w:t></w:r>" +
         "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
         "<w:r><w:t xml:space='preserve'>O SCRUM se adaptaria ao projeto, é um processo ágil com transparência, todos tem conhecimento dos requisitos e os processos e do andamento do processo, </w:t></w:r>" +
         "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" +
         "<w:bookmarkEnd w:id='0'/>" +
         "<w:r><w:t>processo com constante inspeção, seja no sprint review ou nas reuniões diárias além de ser um processo extremamente adaptativo, utilizando este processo o cliente vai poder participar de perto do desenvolvimento, além de ter um rápido resultado, vai ajudar o cliente a decidir em qual ponto focar durante o desenvolvimento, já que ele tem a necessidade de muitos recursos.</w:t></w:r>" +
         "</w:p>"

$blankXml = "<w:p $wns/>"

$insertPoint.InsertXML($blankXml + $p3Xml + $blankXml + $blankXml)

Write-Host 'Final paragraph count:' $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host "Para $i => " $p.Range.Text
}
